# Update COVID-19 UK data as of 29 April 2020
# Column F holds the raw "confirmed cases" input values (one row per
# local authority); columns C/D are formulas derived from F and
# recompute automatically. Columns L/M (and K23/K24) hold the raw
# "deaths" input values for the second table; columns I/J are formulas
# derived from L and recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F: raw case counts per local authority (row -> new value) ---
$fValues = @{ 2 = 448; 3 = 1176; 4 = 608; 7 = 597; 8 = 2782; 9 = 301; 10 = 413; 11 = 732; 12 = 400; 13 = 205; 14 = 796; 15 = 1330; 16 = 354; 17 = 593; 18 = 1027; 19 = 789; 20 = 434; 22 = 718; 23 = 567; 24 = 457; 25 = 813; 26 = 746; 28 = 482; 29 = 1279; 30 = 649; 31 = 1336; 32 = 1793; 33 = 253; 34 = 466; 35 = 1195; 36 = 676; 37 = 426; 38 = 301; 39 = 720; 40 = 1030; 41 = 539; 42 = 552; 43 = 824; 44 = 2451; 46 = 1152; 47 = 620; 49 = 309; 50 = 605; 51 = 2689; 52 = 539; 53 = 886; 54 = 167; 55 = 623; 56 = 282; 57 = 2202; 58 = 711; 59 = 635; 61 = 406; 62 = 456; 63 = 3272; 64 = 422; 65 = 435; 66 = 533; 67 = 516; 68 = 1117; 69 = 2695; 70 = 1322; 71 = 630; 72 = 884; 73 = 898; 75 = 1363; 76 = 575; 77 = 1047; 79 = 637; 81 = 443; 82 = 852; 83 = 938; 84 = 1427; 85 = 130; 86 = 329; 87 = 218; 89 = 872; 90 = 1013; 91 = 732; 93 = 1030; 94 = 582; 95 = 1515; 96 = 223; 98 = 281; 99 = 454; 100 = 658; 102 = 358; 103 = 417; 104 = 666; 106 = 583; 107 = 940; 108 = 770; 109 = 2039; 110 = 464; 111 = 370; 112 = 544; 113 = 464; 114 = 339; 116 = 444; 117 = 290; 118 = 1167; 119 = 614; 120 = 1650; 121 = 796; 122 = 334; 123 = 424; 124 = 999; 126 = 2498; 127 = 656; 128 = 350; 129 = 444; 130 = 244; 132 = 200; 134 = 574; 135 = 494; 136 = 920; 137 = 712; 138 = 904; 139 = 578; 140 = 1116; 141 = 312; 142 = 1083; 143 = 588; 144 = 615; 145 = 405; 146 = 240; 147 = 971; 148 = 362; 149 = 777; 150 = 1081 }
foreach ($row in $fValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $fValues[$row]
}

# --- Column L: raw death counts per region (row -> new value) ---
$lValues = @{ 2 = 6530; 3 = 10035; 4 = 24090; 5 = 7224; 6 = 18106; 7 = 16559; 8 = 6056; 9 = 12593; 10 = 9665; 14 = 115859; 15 = 3463; 16 = 11034; 17 = 9629 }
foreach ($row in $lValues.Keys) {
    $ws.Cells.Item($row, 12).Value = $lValues[$row]
}

# --- Column M: raw supplementary counts (row -> new value) ---
$mValues = @{ 14 = 23550; 15 = 329; 16 = 1332; 17 = 886 }
foreach ($row in $mValues.Keys) {
    $ws.Cells.Item($row, 13).Value = $mValues[$row]
}

# --- Column K: raw totals (row -> new value) for rows 23/24 ---
$kValues = @{ 23 = 165221; 24 = 26097 }
foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 11).Value = $kValues[$row]
}

# --- Update the active selection to match the saved workbook state ---
$ws.Range("I20").Select() | Out-Null
